# :construction: Red wines added to varieties
$wb = $excel.ActiveWorkbook

# --- "Wine & food pairing": add the Corvina veronese pairing row ---
# (written first so the new shared string for this value keeps the same
#  table position it has in the authored workbook)
$ws7 = $wb.Worksheets.Item("Wine & food pairing")
$ws7.Activate()
$ws7.Range("B10").Value = 'Hare civet, Truffle risotto'
$ws7.Range("A10").Value = "Corvina veronese"
$ws7.Range("B11").Select()

# --- "Tastes & aromas of varieties": fill in Red-wine taste/aroma descriptions (col B) ---
$ws6 = $wb.Worksheets.Item("Tastes & aromas of varieties")
$ws6.Activate()

$ws6.Range("B23").Value = 'Strong tanin taste, strong to medium strenght, cassis aromas, and sometime vegetal aromas, or sometime other fruits; Can age (Red wine)'
$ws6.Range("B24").Value = 'Deep color, deep body, strong level of alcohol and low tanin taste; Plum or chocolate aromas;(Red wine)'
$ws6.Range("B25").Value = 'Strong tanin and acidity; Strong alcohol level that makes it smoother; Dark when it''s young, and starts to get some orange shading when it ages; Fruits aromas: strawberry, jam; Earthy and wooden aromas: tar, 0 (Red wine)'
$ws6.Range("B26").Value = 'Delicate and stimulating; Paler than Merlot or Cabernet; Low to medium acidity; Fruity aromas: Red berries and earthy or wooden aromas (Red wine)'
$ws6.Range("B27").Value = 'Low to medium acidity; Medium tanin; Fruity aromas, cherry; Floral aromas, violet; Sometime hazelnut aromas (Red wine)'
$ws6.Range("B28").Value = 'Dark color; Very robust; Strong tanin level; Aromas: smoked meat, grilled pepper, tar, burnt rubber, sometime Strawberry; (Red wine)'
$ws6.Range("B29").Value = 'Rich wine; High level of alcohol; Medium to high level of tanin; Aromas: blackberry, raspberry, spicy notes, jam (Red wine)'
$ws6.Range("B30").Value = 'Powerful wines; Can age well (Red wine)'
$ws6.Range("B31").Value = 'Light color; Low tanin; Fruity aromas (Red wine)'
$ws6.Range("B32").Value = 'Light wines; Delicates, fruity; Good acidity and low tanin; (Red wine)'
$ws6.Range("B33").Value = 'Elegant tanin structure, nice aromas, very fruity; good ageing capacity (Red wine)'
$ws6.Range("B34").Value = 'Powerful aromas and spicy; Sweet (Red wine)'
$ws6.Range("B35").Value = 'Powerful and rich with great fineness (Red wine)'
$ws6.Range("B36").Value = 'Fruits and violet aromas (Red wine)'
$ws6.Range("B37").Value = 'Spicy and fruity aromas; Good ageing capacity (Red wine)'
$ws6.Range("B38").Value = 'Subtil and delicate; Similar to cabernet franc (Red wine)'
$ws6.Range("B39").Value = 'Strong, structured, robust, pulpy fruits, rich, sappy with bitter almond notes  (Red wine)'
$ws6.Range("B40").Value = 'Kirsch aromas (Red wine)'

$ws6.Range("B41").Select()
